$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 683, pushing existing rows 683:771 down to 684:772.
$ws.Rows.Item(683).Insert()

# Populate the newly inserted row 683 with the new weekly record.
$ws.Cells.Item(683, 1).Value  = 10
$ws.Cells.Item(683, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(683, 3).Value  = "La Araucanía"
$ws.Cells.Item(683, 4).Value  = 45131
$ws.Cells.Item(683, 4).NumberFormat = $ws.Cells.Item(684, 4).NumberFormat
$ws.Cells.Item(683, 5).Value  = 9
$ws.Cells.Item(683, 6).Value  = "Fruta"
$ws.Cells.Item(683, 7).Value  = 100108
$ws.Cells.Item(683, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(683, 9).Value  = 100108005
$ws.Cells.Item(683, 10).Value = "Piña"
$ws.Cells.Item(683, 11).Value = "Caramelo"
$ws.Cells.Item(683, 12).Value = "Segunda"
$ws.Cells.Item(683, 13).Value = 1100
$ws.Cells.Item(683, 14).Value = 13000
$ws.Cells.Item(683, 15).Value = 14000
$ws.Cells.Item(683, 16).Value = 13091
$ws.Cells.Item(683, 17).Value = "$/caja 7 unidades"
$ws.Cells.Item(683, 18).Value = "Ecuador"
$ws.Cells.Item(683, 19).Value = 1870
$ws.Cells.Item(683, 20).Value = 7
